$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170..293 down to 171..294
$ws.Rows.Item(170).Insert()

# Populate the new row 170 with the new record's data
$ws.Cells.Item(170, 1).Value = 10
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170, 3).Value = "La Araucanía"
$ws.Cells.Item(170, 4).Value = 44574
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 100112040
$ws.Cells.Item(170, 7).Value = "Cilantro"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 235
$ws.Cells.Item(170, 11).Value = 8000
$ws.Cells.Item(170, 12).Value = 9000
$ws.Cells.Item(170, 13).Value = 8532
$ws.Cells.Item(170, 14).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(170, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(170, 16).Value = 4266
$ws.Cells.Item(170, 17).Value = 2
$ws.Cells.Item(170, 18).Value = "Hortaliza"
